$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.628.47"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "3.678.48"
$ws.Range("E3").Value = "  -0.34%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'646.40"
$ws.Range("E5").Value = "  -4.93%  "
$ws.Range("D6").Value = "'159.78"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'0.502"
$ws.Range("E8").Value = "  +1.35%  "
$ws.Range("E9").Value = "  -1.09%  "
$ws.Range("D10").Value = "'7.17"
$ws.Range("E10").Value = "  -0.13%  "
$ws.Range("E11").Value = "  +0.60%  "
$ws.Range("D12").Value = "'0.0000231"
$ws.Range("E12").Value = "  -1.09%  "
$ws.Range("D13").Value = "4.298.35"
$ws.Range("D14").Value = "'32.67"
$ws.Range("E14").Value = "  +0.53%  "
$ws.Range("D15").Value = "3.670.91"
$ws.Range("E15").Value = "  -0.58%  "
$ws.Range("D16").Value = "69.592.40"
$ws.Range("E16").Value = "  +0.23%  "
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("D18").Value = "'16.01"
$ws.Range("E18").Value = "  -0.30%  "
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("D20").Value = "'468.79"
$ws.Range("E20").Value = "  -0.56%  "
$ws.Range("D21").Value = "'10.10"
$ws.Range("E21").Value = "  +2.96%  "
$ws.Range("E22").Value = "  -0.25%  "
$ws.Range("E23").Value = "  -0.95%  "
$ws.Range("D24").Value = "3.824.04"
$ws.Range("E24").Value = "  -0.33%  "
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("E26").Value = "  -1.55%  "
$ws.Range("D27").Value = "'10.98"
$ws.Range("E27").Value = "  +0.85%  "
$ws.Range("D28").Value = "'9.09"
$ws.Range("E28").Value = "  -0.85%  "
$ws.Range("E29").Value = "  -2.43%  "
$ws.Range("E30").Value = "  -1.81%  "
$ws.Range("D31").Value = "'2.00"
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("E32").Value = "  -0.22%  "
$ws.Range("E33").Value = "  -0.89%  "
$ws.Range("E34").Value = "  -2.10%  "
$ws.Range("D35").Value = "3.675.21"
$ws.Range("E35").Value = "  -0.14%  "
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("D37").Value = "'8.44"
$ws.Range("E37").Value = "  -0.27%  "
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("E39").Value = "  -5.20%  "
$ws.Range("D40").Value = "'1.00"
$ws.Range("D41").Value = "'176.26"
$ws.Range("E41").Value = "  +4.61%  "
$ws.Range("D42").Value = "'2.21"
$ws.Range("E42").Value = "  -2.19%  "
$ws.Range("D43").Value = "'0.0900"
$ws.Range("E43").Value = "  -0.45%  "
$ws.Range("D44").Value = "'0.926"
$ws.Range("E44").Value = "  -1.60%  "
$ws.Range("D45").Value = "'47.22"
$ws.Range("E45").Value = "  +1.00%  "
$ws.Range("D46").Value = "'28.87"
$ws.Range("E46").Value = "  +3.13%  "
$ws.Range("E47").Value = "  -1.73%  "
$ws.Range("E48").Value = "  -1.52%  "
$ws.Range("E49").Value = "  -5.59%  "
$ws.Range("D50").Value = "'7.82"
$ws.Range("E50").Value = "  -0.87%  "
$ws.Range("E51").Value = "  -3.77%  "
